# Applies:
#  1. "Youtube tutorials:" -> split "Youtube" out with spellStart/spellEnd proofErr marks
#  2. "25/9/2024 <EnDash> watched figma tutorials on youtube" -> split "figma" and
#     "youtube" out with spellStart/spellEnd proofErr marks
#  3. Append a new run ", create article categories page" right after the run
#     ", created connections between page"

$d = $word.ActiveDocument
$enDash = [char]0x2013

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$pkgClose = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1. "Youtube tutorials:" ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Youtube tutorials:`r") {
        $rng = $p.Range
        $sub = $d.Range($rng.Start, $rng.End - 1)
        $sub.Text = ""
        $inner = '<w:proofErr w:type="spellStart"/>' +
                 '<w:r><w:t>Youtube</w:t></w:r>' +
                 '<w:proofErr w:type="spellEnd"/>' +
                 '<w:r><w:t xml:space="preserve"> tutorials:</w:t></w:r>'
        $sub.InsertXML($pkgOpen + $inner + $pkgClose)
        break
    }
}

# --- 2. "25/9/2024 - watched figma tutorials on youtube" ---
$target2 = "25/9/2024 " + $enDash + " watched figma tutorials on youtube"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq ($target2 + "`r")) {
        $rng = $p.Range
        $sub = $d.Range($rng.Start, $rng.End - 1)
        $sub.Text = ""
        $inner = '<w:r><w:t xml:space="preserve">25/9/2024 &#8211; watched </w:t></w:r>' +
                 '<w:proofErr w:type="spellStart"/>' +
                 '<w:r><w:t>figma</w:t></w:r>' +
                 '<w:proofErr w:type="spellEnd"/>' +
                 '<w:r><w:t xml:space="preserve"> tutorials on </w:t></w:r>' +
                 '<w:proofErr w:type="spellStart"/>' +
                 '<w:r><w:t>youtube</w:t></w:r>' +
                 '<w:proofErr w:type="spellEnd"/>'
        $sub.InsertXML($pkgOpen + $inner + $pkgClose)
        break
    }
}

# --- 3. Append new run after ", created connections between page" ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "created connections between page") {
        $p.Range.InsertAfter(", create article categories page")
        break
    }
}
